$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K6").Value = 7503108.600000001
$ws.Range("I6").Value = 2501036.2
$ws.Range("M6").Value = -7502996.600000001
$ws.Range("H6").Value = 1432449
$ws.Range("M9").Value = -6848.6924
$ws.Range("I9").Value = 7017.6924
$ws.Range("H9").Value = 6155
$ws.Range("L9").Value = 547.5
$ws.Range("J9").Value = 547.5
$ws.Range("N9").Value = -885.5
$ws.Range("K9").Value = 7017.6924
$ws.Range("I11").Value = 21.875
$ws.Range("H11").Value = 21.875
$ws.Range("K11").Value = 21.875
$ws.Range("M11").Value = 118.125
$ws.Range("H76").Value = 3749.5
$ws.Range("I76").Value = 4000
$ws.Range("M76").Value = -3685
$ws.Range("K76").Value = 4000
$ws.Range("H79").Value = 3749.5
$ws.Range("I79").Value = 4000
$ws.Range("M79").Value = -2908
$ws.Range("K79").Value = 4000
$ws.Range("J86").Value = 6647.3335
$ws.Range("L86").Value = 6647.3335
$ws.Range("K86").Value = 5142.857
$ws.Range("N86").Value = -8893.333500000001
$ws.Range("H86").Value = 5594.2
$ws.Range("I86").Value = 5142.857
$ws.Range("M86").Value = -4019.857
$ws.Range("I89").Value = 5142.857
$ws.Range("J89").Value = 6647.3335
$ws.Range("H89").Value = 5594.2
$ws.Range("K89").Value = 25714.285
$ws.Range("L89").Value = 33236.6675
$ws.Range("N89").Value = -44468.6675
$ws.Range("M89").Value = -20098.285
$ws.Range("K94").Value = 21939.666
$ws.Range("N94").Value = -30902
$ws.Range("L94").Value = 30000
$ws.Range("I94").Value = 21939.666
$ws.Range("M94").Value = -21488.666
$ws.Range("J94").Value = 30000
$ws.Range("H94").Value = 22887.941
$ws.Range("H98").Value = 743.2143
$ws.Range("I98").Value = 754.96155
$ws.Range("K98").Value = 754.96155
$ws.Range("M98").Value = 743.03845
$ws.Range("M100").Value = -2792.1667
$ws.Range("I100").Value = 3333.1667
$ws.Range("K100").Value = 3333.1667
$ws.Range("H100").Value = 3354.2222
$ws.Range("J101").Value = 274.66666
$ws.Range("N101").Value = -4067.99998
$ws.Range("H101").Value = 327.1
$ws.Range("L101").Value = 823.9999799999999
$ws.Range("H105").Value = 39474
$ws.Range("N105").Value = -46462
$ws.Range("L105").Value = 39474
$ws.Range("J105").Value = 39474
$ws.Range("H122").Value = 743.2143
$ws.Range("I122").Value = 754.96155
$ws.Range("K122").Value = 2264.88465
$ws.Range("M122").Value = 185.11535
$ws.Range("K132").Value = 33438249
$ws.Range("I132").Value = 11146083
$ws.Range("M132").Value = -33435719
$ws.Range("H132").Value = 11146083
$ws.Range("M135").Value = -2661.17673
$ws.Range("I135").Value = 577.35297
$ws.Range("K135").Value = 5196.17673
$ws.Range("H135").Value = 767.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I32").Value = 2035.1803
$ws.Range("K32").Value = 2035.1803
$ws.Range("M32").Value = -1748.1803
$ws.Range("H32").Value = 2285.0938
$ws.Range("M45").Value = -11907.125
$ws.Range("I45").Value = 12284.125
$ws.Range("H45").Value = 14876.667
$ws.Range("K45").Value = 12284.125
$ws.Range("M97").Value = -654.0358000000001
$ws.Range("H97").Value = 1150.0358
$ws.Range("J97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("K97").Value = 1150.0358
$ws.Range("L97").Value = 0
$ws.Range("I97").Value = 1150.0358

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N94").Value = -4735.3333
$ws.Range("K94").Value = 1225.2609
$ws.Range("L94").Value = 3833.3333
$ws.Range("I94").Value = 1225.2609
$ws.Range("M94").Value = -774.2609
$ws.Range("J94").Value = 3833.3333
$ws.Range("H94").Value = 1526.1923
$ws.Range("K99").Value = 1570.2
$ws.Range("M99").Value = -72.20000000000005
$ws.Range("J99").Value = 12556.857
$ws.Range("L99").Value = 12556.857
$ws.Range("N99").Value = -15552.857
$ws.Range("H99").Value = 3973.5312
$ws.Range("I99").Value = 1570.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 25533.809
$ws.Range("K31").Value = 25533.809
$ws.Range("H31").Value = 23313.94
$ws.Range("M31").Value = -25238.809
$ws.Range("K34").Value = 25533.809
$ws.Range("I34").Value = 25533.809
$ws.Range("H34").Value = 23313.94
$ws.Range("M34").Value = -25331.809
$ws.Range("M107").Value = 1305.4286
$ws.Range("K107").Value = 614.5714
$ws.Range("I107").Value = 614.5714
$ws.Range("H107").Value = 603.45715
$ws.Range("N132").Value = -22902.5
$ws.Range("K132").Value = 10467.714
$ws.Range("I132").Value = 3489.238
$ws.Range("L132").Value = 17842.5
$ws.Range("M132").Value = -7937.714
$ws.Range("J132").Value = 5947.5
$ws.Range("H132").Value = 3703

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J3").Value = 17000
$ws.Range("M3").Value = -750004358
$ws.Range("I3").Value = 250001490
$ws.Range("L3").Value = 51000
$ws.Range("N3").Value = -51224
$ws.Range("H3").Value = 200004590
$ws.Range("K3").Value = 750004470
$ws.Range("K8").Value = 3244.3335
$ws.Range("I8").Value = 1081.4445
$ws.Range("H8").Value = 1081.4445
$ws.Range("M8").Value = -3105.3335
$ws.Range("H40").Value = 127.2
$ws.Range("L40").Value = 915
$ws.Range("J40").Value = 228.75
$ws.Range("N40").Value = -1053
$ws.Range("K86").Value = 1338
$ws.Range("H86").Value = 578
$ws.Range("I86").Value = 446
$ws.Range("M86").Value = -152
$ws.Range("I89").Value = 446
$ws.Range("H89").Value = 578
$ws.Range("K89").Value = 4014
$ws.Range("M89").Value = 1914
$ws.Range("N124").Value = -39820
$ws.Range("J124").Value = 10000
$ws.Range("L124").Value = 30000
$ws.Range("H124").Value = 9923.076999999999
$ws.Range("L131").Value = 23481.8775
$ws.Range("H131").Value = 49466.15
$ws.Range("M131").Value = -996955.02
$ws.Range("J131").Value = 7827.2925
$ws.Range("I131").Value = 333998.34
$ws.Range("K131").Value = 1001995.02
$ws.Range("N131").Value = -33561.8775
$ws.Range("I137").Value = 3848.5
$ws.Range("M137").Value = -6445.5
$ws.Range("L137").Value = 25425
$ws.Range("J137").Value = 8475
$ws.Range("N137").Value = -35625
$ws.Range("K137").Value = 11545.5
$ws.Range("H137").Value = 4689.6816

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K132").Value = 633565.3500000001
$ws.Range("I132").Value = 211188.45
$ws.Range("M132").Value = -631035.3500000001
$ws.Range("H132").Value = 202795.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I7").Value = 16506.857
$ws.Range("L7").Value = 9000
$ws.Range("H7").Value = 15568.5
$ws.Range("M7").Value = -16394.857
$ws.Range("N7").Value = -9224
$ws.Range("K7").Value = 16506.857
$ws.Range("J7").Value = 9000
$ws.Range("L93").Value = 9633
$ws.Range("H93").Value = 698763.9399999999
$ws.Range("J93").Value = 9633
$ws.Range("N93").Value = -12129
$ws.Range("M100").Value = -62861.59
$ws.Range("I100").Value = 63402.59
$ws.Range("N100").Value = -14015.333
$ws.Range("J100").Value = 12933.333
$ws.Range("L100").Value = 12933.333
$ws.Range("K100").Value = 63402.59
$ws.Range("H100").Value = 55832.2
$ws.Range("H122").Value = 10491.875
$ws.Range("I122").Value = 7533.727
$ws.Range("K122").Value = 22601.181
$ws.Range("M122").Value = -20151.181
$ws.Range("H126").Value = 15568.5
$ws.Range("N126").Value = -31940
$ws.Range("M126").Value = -47050.571
$ws.Range("K126").Value = 49520.571
$ws.Range("L126").Value = 27000
$ws.Range("I126").Value = 16506.857
$ws.Range("J126").Value = 9000
$ws.Range("K132").Value = 9411.545999999998
$ws.Range("I132").Value = 3137.182
$ws.Range("M132").Value = -6881.545999999998
$ws.Range("H132").Value = 3731.4614

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N15").Value = -23583
$ws.Range("H15").Value = 23007
$ws.Range("J15").Value = 23007
$ws.Range("L15").Value = 23007
$ws.Range("N132").Value = -8057
$ws.Range("K132").Value = 8392.235700000001
$ws.Range("I132").Value = 2797.4119
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = -5862.235700000001
$ws.Range("J132").Value = 999
$ws.Range("H132").Value = 2697.5
